# Append a new data row (row 92) to Sheet1 of the profit-tracking workbook,
# matching the run performed on 2026-02-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/24/2026"
$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).Value = 8814.75
$ws.Cells.Item($row, 3).Value = 0.2451242745093181
$ws.Cells.Item($row, 4).Value = 0.7548757254906819
$ws.Cells.Item($row, 5).Value = -362.64
$ws.Cells.Item($row, 6).Value = -38.79
$ws.Cells.Item($row, 7).Value = -24337.71
$ws.Cells.Item($row, 8).Value = -78.53
$ws.Cells.Item($row, 9).Value = -1292.14
$ws.Cells.Item($row, 10).Value = -37.42
$ws.Cells.Item($row, 11).Value = -25629.85
$ws.Cells.Item($row, 12).Value = -74.41
